$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in values for previously-blank cells in rows 10 and 11 ---
$ws.Range("D10").Value = 23.2
$ws.Range("E10").Value = 22.6
$ws.Range("F10").Value = 21.9
$ws.Range("H10").Value = 30.7
$ws.Range("I10").Value = 35.4
$ws.Range("J10").Value = 28.5

$ws.Range("D11").Value = 41.7
$ws.Range("E11").Value = 36.8
$ws.Range("F11").Value = 37.1
$ws.Range("H11").Value = 40.4
$ws.Range("I11").Value = 39.6
$ws.Range("J11").Value = 41.4

# --- Add a new column T (year 2023) mirroring the formatting of column S ---
# Row -> 2023 value, for every row that already carries a value in column S.
$newColumnValues = [ordered]@{
   4  = 2023
   5  = 29.810232786618478
   7  = 29.669466599025686
   8  = 29.964546620904322
   10 = 30.196132774743152
   11 = 29.586576623908091
   13 = 36.679314997357302
   14 = 27.761651579699627
   15 = 26.945499612171261
   16 = 19.920707357966336
   18 = 48.132487638243802
   19 = 47.095468608697217
   20 = 49.224436679851941
   21 = 36.060409324309092
   22 = 36.438841751655779
   23 = 35.670056408825062
   24 = 30.852411825788565
   25 = 30.330534800771165
   26 = 31.403214482728419
   27 = 38.094443042646382
   28 = 38.158675127516169
   29 = 38.024712001344874
   30 = 20.396452079475392
   31 = 20.655435559889604
   32 = 20.117912106064367
   33 = 23.188885535955222
   34 = 23.078628487453106
   35 = 23.300729383023359
   36 = 26.600808028614065
   37 = 26.440501693813694
   38 = 26.779307280034676
   39 = 32.389629740110649
   40 = 31.71205247152805
   41 = 33.231499210635427
   42 = 13.872134221368512
   43 = 14.562707317462859
   44 = 13.080583219648313
}

# Rows that carry formatting in column S but have no numeric value there -
# column T should pick up the same (blank) formatting.
$blankFormatRows = @(6)

foreach ($r in $newColumnValues.Keys) {
    $srcCell = $ws.Range("S$r")
    $dstCell = $ws.Range("T$r")
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null
    $dstCell.Value = $newColumnValues[$r]
}

foreach ($r in $blankFormatRows) {
    $srcCell = $ws.Range("S$r")
    $dstCell = $ws.Range("T$r")
    $srcCell.Copy() | Out-Null
    $dstCell.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# --- Drop the stale cell selection left over from the previous save ---
$ws.Range("A1").Select() | Out-Null
